# Generate Report for Handoff
# - Updates the "Status" text from "In Translation" to "Ready for handoff"
#   and the associated "Latest HO Xliff Generate Date" / "Latest Handoff
#   Datetime" timestamps on the Overview / zh-cn / de-de sheets.
# - Widens the Status-related columns (col 5/6 on Overview, col 3 on the
#   language sheets) to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: E2/F2 hold the status text, G2 holds the handoff-generate
# timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-12 19:13:21"

# zh-cn / de-de sheets: C2 holds the Status text, H2 holds the
# "Latest Handoff Datetime" timestamp.
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-12 19:13:14"
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-12 19:13:21"

# Widen the Status-adjacent columns to match the new text width. The stored
# OOXML <col width=.../> is ~5/6 of a character wider than the COM
# ColumnWidth value (Excel's standard internal-margin offset for this
# engine's metrics), so back that out to land as close as possible to the
# target stored width of 17.2159881591797.
$newStatusColWidth = 17.2159881591797 - (5 / 6)
$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
